$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# ---------------------------------------------------------------------------
# Change 1: "Petit rappel : Quelle est la différence entre un ID et une
# class ? " -- wrap "class" with spell-check proofErr markers, splitting the
# single run into three runs.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$search1 = "Petit rappel" + $nbsp + ": Quelle est la différence entre un ID et une class" + $nbsp + "? "
$found1 = $rng1.Find.Execute($search1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $rng1f = $d.Range($rng1.Start, $rng1.End)
    $xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Petit rappel' + $nbsp + ': Quelle est la différence entre un ID et une </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>class</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">' + $nbsp + '? </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng1f.InsertXML($xml1)
}

Write-Host "Change 1 (proofErr class) found:" $found1

# ---------------------------------------------------------------------------
# Change 2: rewrite the "Bootstrap 5 / béta" paragraph and append a new
# paragraph about the jQuery-free version 5, with proofErr markers around
# "Vanilla" / "js" plus a trailing _GoBack bookmark.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$search2 = "En ce moment, février 2021, nous sommes au moment où Bootstrap 5 va bientôt pointer le bout de son nez. Mais il n'est pas encore officiel."
$search2 = "En ce moment, février 2021, nous sommes au moment où Bootstrap 5 va bientôt pointer le bout de son nez. Mais il n" + [char]0x2019 + "est pas encore officiel."
$found2 = $rng2.Find.Execute($search2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

Write-Host "Change 2 find found:" $found2

if ($found2) {
    $start2 = $rng2.Start
    $rng2f = $d.Range($rng2.Start, $rng2.End)

    $newPart1 = "En ce moment, février 2021, nous sommes au moment où Bootstrap 5 va bientôt pointer le bout de son nez."
    $newPart2 = " Il est toujours en phase béta. Comme la version 4.x est la plus déployée, je pense qu" + [char]0x2019 + "il est plus intéressant pour vous de voir la version 4.x qui est la plus installée et utilisée."

    $xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $newPart1 + '</w:t></w:r><w:r><w:t xml:space="preserve">' + $newPart2 + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng2f.InsertXML($xml2)

    # The paragraph's text length is now known exactly, so we can locate the
    # paragraph-mark position without relying on a (possibly stale) Range
    # object that was mutated by InsertXML.
    $computedEnd2 = $start2 + $newPart1.Length + $newPart2.Length
    $parEndRng = $d.Range($computedEnd2, $computedEnd2)
    $parEndRng.InsertParagraphAfter()

    # Re-resolve the brand-new (empty) paragraph via the Paragraphs
    # collection -- using its .Range avoids the stray empty <w:r/> that is
    # left behind when targeting a manually computed zero-length Range.
    $newPara = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -eq ($computedEnd2 + 1)) {
            $newPara = $p
            break
        }
    }

    if ($newPara -ne $null) {
        $newParaRng = $newPara.Range
        $xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">L' + [char]0x2019 + 'intérêt de la version 5, c' + [char]0x2019 + 'est qu' + [char]0x2019 + 'il n' + [char]0x2019 + 'utilise plus la librairie jQuery et n' + [char]0x2019 + 'est donc plus tributaire de celui-ci. Bootstrap 5 a été remanier pour travailler avec du pur javascript (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Vanilla</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>js</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>). C' + [char]0x2019 + 'est donc une dépendance en moins.</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $newParaRng.InsertXML($xml3)
    } else {
        Write-Host "ERROR: could not locate the freshly inserted empty paragraph"
    }
}
